$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.404.41"
$ws.Range("E2").Value = "  -1.81%  "

$ws.Range("D3").Value = "1.792.75"
$ws.Range("E3").Value = "  -2.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4561"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.13%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.02%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.277"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.29%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.722.93"
$ws.Range("E15").Value = "  -6.63%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.319"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.94%  "

$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008566"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.12%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.442.10"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.987"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "

$ws.Range("D24").Value = "2.013.35"
$ws.Range("E24").Value = "  -3.33%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.975"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.038"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "112.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.856"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.99%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08656"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.072"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.81%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.445"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7231"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.88%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.663"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.110"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "

$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.006"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.077"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01944"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.876"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.08%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5244"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.20%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.893"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.57%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1526"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.58%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.014"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.63%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4693"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.886"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.25%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.81%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.585"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "

